# Quarterly symbol-list refresh: updated Price (D) and Volume(1h) (E) columns
# for the rows whose source data changed, per the GitHub Actions data pull.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Leading apostrophe forces these to stay text cells (matching the sheet's
# existing text-typed Price/Volume columns) instead of being auto-converted
# to numbers/percentages by Excel; resetting the style afterwards keeps the
# cell formatting identical to the untouched cells around it.
$ws.Range("D2").Value = "'327.22"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'-0.96%"
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "'43.73"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'5.06%"
$ws.Range("E3").Style = "Normal"
$ws.Range("D4").Value = "'5.485"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "'-3.82%"
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = "'0.08057"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'-4.36%"
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Value = "'8.639"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'-2.12%"
$ws.Range("E6").Style = "Normal"
$ws.Range("E7").Value = "'-4.39%"
$ws.Range("E7").Style = "Normal"
$ws.Range("D8").Value = "'1.876"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "'-5.80%"
$ws.Range("E8").Style = "Normal"
$ws.Range("D9").Value = "'2.698"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'-7.71%"
$ws.Range("E9").Style = "Normal"
$ws.Range("D10").Value = "'0.9370"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'1.12%"
$ws.Range("E10").Style = "Normal"
$ws.Range("D11").Value = "'0.1188"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'-7.07%"
$ws.Range("E11").Style = "Normal"
$ws.Range("D12").Value = "'0.1892"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'-3.68%"
$ws.Range("E12").Style = "Normal"
$ws.Range("D13").Value = "'0.09480"
$ws.Range("D13").Style = "Normal"
$ws.Range("D14").Value = "'0.04162"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'5.34%"
$ws.Range("E14").Style = "Normal"
$ws.Range("E15").Value = "'0.36%"
$ws.Range("E15").Style = "Normal"
$ws.Range("D16").Value = "'0.001266"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'-2.95%"
$ws.Range("E16").Style = "Normal"
$ws.Range("D17").Value = "'0.006040"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'-1.24%"
$ws.Range("E17").Style = "Normal"
$ws.Range("E18").Value = "'4.24%"
$ws.Range("E18").Style = "Normal"
$ws.Range("D19").Value = "'0.3485"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'-0.41%"
$ws.Range("E19").Style = "Normal"
$ws.Range("D20").Value = "'8.583"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "'-4.16%"
$ws.Range("E20").Style = "Normal"
$ws.Range("E21").Value = "'0.05%"
$ws.Range("E21").Style = "Normal"
$ws.Range("D22").Value = "'0.2594"
$ws.Range("D22").Style = "Normal"
$ws.Range("D23").Value = "'0.04341"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'-1.85%"
$ws.Range("E23").Style = "Normal"
$ws.Range("E24").Value = "'-0.95%"
$ws.Range("E24").Style = "Normal"
$ws.Range("D25").Value = "'0.004348"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'-0.61%"
$ws.Range("E25").Style = "Normal"
$ws.Range("D26").Value = "'0.0001233"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "'3.50%"
$ws.Range("E26").Style = "Normal"
$ws.Range("D27").Value = "'0.0004001"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "'0.09%"
$ws.Range("E27").Style = "Normal"
$ws.Range("D39").Value = "'0.02637"
$ws.Range("D39").Style = "Normal"
$ws.Range("D40").Value = "'0.05380"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'-2.53%"
$ws.Range("E40").Style = "Normal"
$ws.Range("D41").Value = "'0.01097"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'22.25%"
$ws.Range("E41").Style = "Normal"
$ws.Range("D42").Value = "'0.007628"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'-3.57%"
$ws.Range("E42").Style = "Normal"
$ws.Range("D43").Value = "'0.1386"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'-3.58%"
$ws.Range("E43").Style = "Normal"
$ws.Range("D44").Value = "'0.002128"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'1.69%"
$ws.Range("E44").Style = "Normal"
$ws.Range("D45").Value = "'0.009699"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'-11.86%"
$ws.Range("E45").Style = "Normal"
$ws.Range("D46").Value = "'0.00006875"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'-6.39%"
$ws.Range("E46").Style = "Normal"
$ws.Range("E47").Value = "'0.10%"
$ws.Range("E47").Style = "Normal"
$ws.Range("D48").Value = "'0.003570"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "'9.88%"
$ws.Range("E48").Style = "Normal"
$ws.Range("D49").Value = "'0.002277"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "'-0.24%"
$ws.Range("E49").Style = "Normal"
$ws.Range("E50").Value = "'0.10%"
$ws.Range("E50").Style = "Normal"
$ws.Range("D51").Value = "'0.0002006"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "'0.10%"
$ws.Range("E51").Style = "Normal"
